$d = $word.ActiveDocument

# 1. Title: "Debate 1" -> "Debate 4"
$d.Content.Find.Execute("Debate 1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Debate 4", 2) | Out-Null

# 2. Date: "F2025" -> "S2026"
$d.Content.Find.Execute("F2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "S2026", 2) | Out-Null

# Locate the "Date" paragraph (so we don't depend on a hard-coded index).
$dateIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Style.NameLocal -eq "Date") {
        $dateIdx = $i
        break
    }
}

# 3. Insert a new "Overview" Heading2 paragraph right after the Date paragraph.
$datePara = $d.Paragraphs($dateIdx)
$datePara.Range.InsertParagraphAfter()
$overviewPara = $d.Paragraphs($dateIdx + 1)
$overviewPara.Style = "Heading2"
$overviewPara.Range.Text = "Overview"

# 4. FirstParagraph body text: "Debate stuff" -> "Nothing to see here yet"
$d.Content.Find.Execute("Debate stuff", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Nothing to see here yet", 2) | Out-Null

# 5. Wrap the new "Overview" section -- the heading paragraph together with
#    the body paragraph right after it -- in a bookmark named "overview",
#    mirroring the bookmarkStart/bookmarkEnd pair added around that section.
$startPara = $d.Paragraphs($dateIdx + 1)
$endPara = $d.Paragraphs($dateIdx + 2)
$bmRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$d.Bookmarks.Add("overview", $bmRange) | Out-Null

Write-Output "done"
